# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-19, replacing the previous Strike#-derived values
$kValues = @{
    2  = 6
    3  = 7
    4  = 0
    5  = 3
    6  = 0
    7  = 5
    8  = 1
    9  = 0
    10 = 2
    11 = 5
    12 = 5
    13 = 2
    14 = 3
    15 = 2
    16 = 3
    17 = 1
    18 = 0
    19 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
